$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap province names whose shared-string rows were reordered in the source data
$tmp = $ws.Range("A15").Value2
$ws.Range("A15").Value2 = $ws.Range("A16").Value2
$ws.Range("A16").Value2 = $tmp

$tmp = $ws.Range("A35").Value2
$ws.Range("A35").Value2 = $ws.Range("A36").Value2
$ws.Range("A36").Value2 = $tmp

$tmp = $ws.Range("A41").Value2
$ws.Range("A41").Value2 = $ws.Range("A42").Value2
$ws.Range("A42").Value2 = $tmp

$tmp = $ws.Range("A47").Value2
$ws.Range("A47").Value2 = $ws.Range("A48").Value2
$ws.Range("A48").Value2 = $tmp

# Update the "Datos actualizados" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 12:52"

# Refresh the numeric data (Casos totales / Casos activos / Recuperados / Muertes)
$ws.Range("B8").Value = 5138
$ws.Range("C8").Value = 2205
$ws.Range("D8").Value = 9768
$ws.Range("E8").Value = 503
$ws.Range("B9").Value = 4526
$ws.Range("C9").Value = 1301
$ws.Range("D9").Value = 2833
$ws.Range("E9").Value = 392
$ws.Range("B11").Value = 3404
$ws.Range("C11").Value = 2205
$ws.Range("D11").Value = 9768
$ws.Range("E11").Value = 309
$ws.Range("B13").Value = 2993
$ws.Range("C13").Value = 764
$ws.Range("D13").Value = 1902
$ws.Range("E13").Value = 327
$ws.Range("B14").Value = 2984
$ws.Range("C14").Value = 2205
$ws.Range("D14").Value = 9768
$ws.Range("E14").Value = 403
$ws.Range("B15").Value = 2976
$ws.Range("C15").Value = 668
$ws.Range("D15").Value = 1993
$ws.Range("E15").Value = 315
$ws.Range("B16").Value = 2908
$ws.Range("C16").Value = 4680
$ws.Range("D16").Value = 4694
$ws.Range("E16").Value = 245
$ws.Range("B35").Value = 1056
$ws.Range("C35").Value = 2205
$ws.Range("D35").Value = 9768
$ws.Range("E35").Value = 139
$ws.Range("B36").Value = 1055
$ws.Range("C36").Value = 111
$ws.Range("D36").Value = 852
$ws.Range("E36").Value = 92
$ws.Range("B37").Value = 1054
$ws.Range("C37").Value = 224
$ws.Range("D37").Value = 731
$ws.Range("B41").Value = 874
$ws.Range("C41").Value = 2205
$ws.Range("D41").Value = 9768
$ws.Range("E41").Value = 129
$ws.Range("B42").Value = 859
$ws.Range("C42").Value = 325
$ws.Range("D42").Value = 441
$ws.Range("E42").Value = 93
$ws.Range("B47").Value = 491
$ws.Range("D47").Value = 328
$ws.Range("E47").Value = 66
$ws.Range("B48").Value = 480
$ws.Range("C48").Value = 97
$ws.Range("D48").Value = 356
$ws.Range("E48").Value = 27
$ws.Range("B49").Value = 459
$ws.Range("C49").Value = 113
$ws.Range("D49").Value = 303
$ws.Range("E49").Value = 43
